# Applies the edits described by the commit "Se inicio el modulo producto, por concluir":
#   1. The fixed Header & Footer date shown on the slide master and every slide layout
#      changes from 19/6/2021 to 25/6/2021.
#   2. The placeholder word "Persona" on slide 5 (inside the rotated "Rectángulo 4"
#      shape, itself nested in the "Grupo 5" group) becomes "producto".

$p = $ppt.ActivePresentation

$oldDate = "19/6/2021"
$newDate = "25/6/2021"

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master footer date placeholder.
$master = $p.SlideMaster
Update-DateShapes $master.Shapes

# Every layout under the master has its own (inherited) copy of the date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes
}

# Slide 5: "Grupo 5" group -> "Rectángulo 4" shape holding the word "Persona".
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $top = $slide5.Shapes.Item($i)
    if ($top.Name -eq "Grupo 5") {
        for ($j = 1; $j -le $top.GroupItems.Count; $j++) {
            $item = $top.GroupItems.Item($j)
            if ($item.HasTextFrame -and $item.TextFrame.TextRange.Text -eq "Persona") {
                $item.TextFrame.TextRange.Text = "producto"
            }
        }
    }
}
